# Missing rows logic and tests
#
# Rename the two existing sheets and add a third sheet ("Missing Rows")
# that reproduces the same 3x3 grid but omits row 2 entirely, to cover
# the "missing rows" code path (as opposed to the existing "Missing
# Cells" sheet, which omits individual cells but keeps every row).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Regular Data"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Missing Cells"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Missing Rows"

$ws3.Range("A1").Value = "A1"
$ws3.Range("B1").Value = "B1"
$ws3.Range("C1").Value = "C1"

# Row 2 intentionally left blank/absent.

$ws3.Range("A3").Value = "A3"
$ws3.Range("B3").Value = "B3"
$ws3.Range("C3").Value = "C3"

$null = $ws3.Range("C3").Select()
